$wb = $excel.ActiveWorkbook

# ===== Sheet: ALC =====
$ws = $wb.Worksheets.Item("ALC")

# Row 40
$ws.Range("H40").Value = 2526.647
$ws.Range("I40").Value = 2395
$ws.Range("J40").Value = 2714.7144
$ws.Range("K40").Value = 2395
$ws.Range("L40").Value = 2714.7144
$ws.Range("M40").Value = -2220
$ws.Range("N40").Value = -3064.7144

# Row 43
$ws.Range("H43").Value = 680
$ws.Range("I43").Value = 625
$ws.Range("K43").Value = 625
$ws.Range("M43").Value = -556

# Row 45
$ws.Range("H45").Value = 1000
$ws.Range("I45").Value = 1000
$ws.Range("K45").Value = 3000
$ws.Range("M45").Value = -2808

# Row 98
$ws.Range("H98").Value = 250039.25
$ws.Range("I98").Value = 386793.75
$ws.Range("J98").Value = 2171.6875
$ws.Range("K98").Value = 386793.75
$ws.Range("L98").Value = 2171.6875
$ws.Range("M98").Value = -385295.75
$ws.Range("N98").Value = -5167.6875

# Row 99
$ws.Range("H99").Value = 221.33333
$ws.Range("I99").Value = 221.33333
$ws.Range("K99").Value = 663.99999
$ws.Range("M99").Value = 834.00001

# Row 112
$ws.Range("H112").Value = 5348479
$ws.Range("J112").Value = 5682715
$ws.Range("L112").Value = 17048145
$ws.Range("N112").Value = -17050361

# Row 122
$ws.Range("H122").Value = 250039.25
$ws.Range("I122").Value = 386793.75
$ws.Range("J122").Value = 2171.6875
$ws.Range("K122").Value = 1160381.25
$ws.Range("L122").Value = 6515.0625
$ws.Range("M122").Value = -1157931.25
$ws.Range("N122").Value = -11415.0625

# Row 129
$ws.Range("H129").Value = 1512.1538
$ws.Range("I129").Value = 444.66666
$ws.Range("J129").Value = 2427.1428
$ws.Range("K129").Value = 1333.99998
$ws.Range("L129").Value = 7281.428400000001
$ws.Range("M129").Value = 3666.00002
$ws.Range("N129").Value = -17281.4284

# Row 137
$ws.Range("H137").Value = 55556852
$ws.Range("I137").Value = 58824740
$ws.Range("K137").Value = 176474220
$ws.Range("M137").Value = -176471670

# Row 138
$ws.Range("H138").Value = 4414531
$ws.Range("I138").Value = 1227711.2
$ws.Range("J138").Value = 6539077.5
$ws.Range("K138").Value = 3683133.6
$ws.Range("L138").Value = 19617232.5
$ws.Range("M138").Value = -3677993.6
$ws.Range("N138").Value = -19627512.5

# ===== Sheet: ARM =====
$ws = $wb.Worksheets.Item("ARM")

# Row 6
$ws.Range("H6").Value = 11500
$ws.Range("I6").Value = 20000
$ws.Range("K6").Value = 20000
$ws.Range("M6").Value = -19827

# Row 32
$ws.Range("H32").Value = 15713.768
$ws.Range("I32").Value = 2117.2188
$ws.Range("J32").Value = 112400.336
$ws.Range("K32").Value = 2117.2188
$ws.Range("L32").Value = 112400.336
$ws.Range("M32").Value = -1830.2188
$ws.Range("N32").Value = -112974.336

# Row 38
$ws.Range("H38").Value = 15002.571
$ws.Range("I38").Value = 15004.5
$ws.Range("J38").Value = 15000
$ws.Range("K38").Value = 15004.5
$ws.Range("L38").Value = 15000
$ws.Range("M38").Value = -14537.5
$ws.Range("N38").Value = -15934

# Row 70
$ws.Range("H70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("N70").ClearContents()

# Row 73
$ws.Range("H73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("N73").ClearContents()

# Row 95
$ws.Range("H95").Value = 500208
$ws.Range("J95").Value = 500208
$ws.Range("L95").Value = 500208
$ws.Range("N95").Value = -505700

# ===== Sheet: BSM =====
$ws = $wb.Worksheets.Item("BSM")

# Row 38
$ws.Range("H38").Value = 30012
$ws.Range("J38").Value = 20018
$ws.Range("L38").Value = 20018
$ws.Range("N38").Value = -20850

# Row 86
$ws.Range("H86").Value = 9193.429
$ws.Range("I86").Value = 3060.2
$ws.Range("J86").Value = 12600.777
$ws.Range("K86").Value = 3060.2
$ws.Range("L86").Value = 12600.777
$ws.Range("M86").Value = -1937.2
$ws.Range("N86").Value = -14846.777

# Row 89
$ws.Range("H89").Value = 9193.429
$ws.Range("I89").Value = 3060.2
$ws.Range("J89").Value = 12600.777
$ws.Range("K89").Value = 15301
$ws.Range("L89").Value = 63003.885
$ws.Range("M89").Value = -9685
$ws.Range("N89").Value = -74235.88500000001

# Row 105
$ws.Range("H105").Value = 3061.743
$ws.Range("I105").Value = 2888.5
$ws.Range("J105").Value = 3354.923
$ws.Range("K105").Value = 2888.5
$ws.Range("L105").Value = 3354.923
$ws.Range("M105").Value = -1141.5
$ws.Range("N105").Value = -6848.923

# Row 134
$ws.Range("H134").Value = 13515416
$ws.Range("I134").Value = 20001250
$ws.Range("J134").Value = 3261.7917
$ws.Range("K134").Value = 60003750
$ws.Range("L134").Value = 9785.375100000001
$ws.Range("M134").Value = -60001215
$ws.Range("N134").Value = -14855.3751

# ===== Sheet: CRP =====
$ws = $wb.Worksheets.Item("CRP")

# Row 22
$ws.Range("H22").Value = 972
$ws.Range("I22").Value = 492.57144
$ws.Range("J22").Value = 2650
$ws.Range("K22").Value = 492.57144
$ws.Range("L22").Value = 2650
$ws.Range("M22").Value = -142.57144
$ws.Range("N22").Value = -3350

# Row 86
$ws.Range("H86").Value = 29413794
$ws.Range("I86").Value = 45456548
$ws.Range("J86").Value = 2083.3333
$ws.Range("K86").Value = 45456548
$ws.Range("L86").Value = 2083.3333
$ws.Range("M86").Value = -45455425
$ws.Range("N86").Value = -4329.3333

# Row 89
$ws.Range("H89").Value = 29413794
$ws.Range("I89").Value = 45456548
$ws.Range("J89").Value = 2083.3333
$ws.Range("K89").Value = 227282740
$ws.Range("L89").Value = 10416.6665
$ws.Range("M89").Value = -227277124
$ws.Range("N89").Value = -21648.6665

# Row 99
$ws.Range("H99").Value = 12501185
$ws.Range("I99").Value = 20834338
$ws.Range("J99").Value = 1457
$ws.Range("K99").Value = 20834338
$ws.Range("L99").Value = 1457
$ws.Range("M99").Value = -20832840
$ws.Range("N99").Value = -4453

# Row 126
$ws.Range("H126").Value = 12501185
$ws.Range("I126").Value = 20834338
$ws.Range("J126").Value = 1457
$ws.Range("K126").Value = 62503014
$ws.Range("L126").Value = 4371
$ws.Range("M126").Value = -62500544
$ws.Range("N126").Value = -9311

# ===== Sheet: CUL =====
$ws = $wb.Worksheets.Item("CUL")

# Row 4
$ws.Range("H4").Value = 25164.285
$ws.Range("I4").Value = 130
$ws.Range("J4").Value = 87750
$ws.Range("K4").Value = 390
$ws.Range("L4").Value = 263250
$ws.Range("M4").Value = -278
$ws.Range("N4").Value = -263474

# Row 63
$ws.Range("H63").Value = 4062.4
$ws.Range("I63").Value = 4062.4
$ws.Range("J63").Value = 0
$ws.Range("K63").Value = 12187.2
$ws.Range("L63").Value = 0
$ws.Range("M63").Value = -11438.2
$ws.Range("N63").ClearContents()

# Row 66
$ws.Range("H66").Value = 4062.4
$ws.Range("I66").Value = 4062.4
$ws.Range("J66").Value = 0
$ws.Range("K66").Value = 36561.6
$ws.Range("L66").Value = 0
$ws.Range("M66").Value = -32817.6
$ws.Range("N66").ClearContents()

# Row 113
$ws.Range("H113").Value = 13889781
$ws.Range("J113").Value = 23810602
$ws.Range("L113").Value = 71431806
$ws.Range("N113").Value = -71436146

# Row 122
$ws.Range("H122").Value = 635.0526
$ws.Range("J122").Value = 1129.125
$ws.Range("L122").Value = 10162.125
$ws.Range("N122").Value = -15062.125

# ===== Sheet: GSM =====
$ws = $wb.Worksheets.Item("GSM")

# Row 11
$ws.Range("H11").Value = 2347208.2
$ws.Range("I11").Value = 3050500.2
$ws.Range("J11").Value = 2901.3333
$ws.Range("K11").Value = 3050500.2
$ws.Range("L11").Value = 2901.3333
$ws.Range("M11").Value = -3050361.2
$ws.Range("N11").Value = -3179.3333

# Row 64
$ws.Range("H64").Value = 20000
$ws.Range("J64").Value = 20000
$ws.Range("L64").Value = 20000
$ws.Range("N64").Value = -20496

# Row 67
$ws.Range("H67").Value = 20000
$ws.Range("J67").Value = 20000
$ws.Range("L67").Value = 20000
$ws.Range("N67").Value = -21716

# Row 132
$ws.Range("H132").Value = 2509.5833
$ws.Range("I132").Value = 2167.6829
$ws.Range("J132").Value = 4512.143
$ws.Range("K132").Value = 6503.048699999999
$ws.Range("L132").Value = 13536.429
$ws.Range("M132").Value = -3973.048699999999
$ws.Range("N132").Value = -18596.429

# ===== Sheet: LTW =====
$ws = $wb.Worksheets.Item("LTW")

# Row 16
$ws.Range("H16").Value = 1453.5454
$ws.Range("I16").Value = 1811.125
$ws.Range("K16").Value = 1811.125
$ws.Range("M16").Value = -1641.125

# Row 46
$ws.Range("H46").Value = 2477.7778
$ws.Range("I46").Value = 1300
$ws.Range("J46").Value = 3066.6667
$ws.Range("K46").Value = 1300
$ws.Range("L46").Value = 3066.6667
$ws.Range("M46").Value = -1112
$ws.Range("N46").Value = -3442.6667

# Row 55
$ws.Range("H55").Value = 468
$ws.Range("I55").Value = 466.33334
$ws.Range("J55").Value = 468.55554
$ws.Range("K55").Value = 466.33334
$ws.Range("L55").Value = 468.55554
$ws.Range("M55").Value = -293.33334
$ws.Range("N55").Value = -814.5555400000001

# Row 87
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

# Row 90
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

# ===== Sheet: WVR =====
$ws = $wb.Worksheets.Item("WVR")

# Row 136
$ws.Range("H136").Value = 9288516
$ws.Range("J136").Value = 4223.6665
$ws.Range("L136").Value = 12670.9995
$ws.Range("N136").Value = -17770.9995
